$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows above the existing table/header so the sheet gains a
# "cover" area with Company Name / Your Name / Phone Number / Email ID
# fields (the "finish attachments" info block).
$ws.Rows("1:4").Insert()

# The table used to live at A1:O2; after the insert it needs to be moved
# back down onto the data that shifted to A5:O6 (Excel does not do this
# automatically because the new rows were inserted above the table, not
# inside it).
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A5:O6"))

# Fill in the new label cells.
$ws.Range("A1").Value = "Company Name:"
$ws.Range("A2").Value = "Your Name:"
$ws.Range("A3").Value = "Phone Number:"
$ws.Range("A4").Value = "Email ID:"

# Style the new rows: labels in column A are bold, and each row alternates
# which of the two columns gets the yellow highlight vs. the accent
# (orange) highlight.
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").Interior.Color = 65535
$ws.Range("B1").Interior.ThemeColor = 6

$ws.Range("A2").Font.Bold = $true
$ws.Range("A2").Interior.ThemeColor = 6
$ws.Range("B2").Interior.Color = 65535

$ws.Range("A3").Font.Bold = $true
$ws.Range("A3").Interior.Color = 65535
$ws.Range("B3").Interior.ThemeColor = 6

$ws.Range("A4").Font.Bold = $true
$ws.Range("A4").Interior.ThemeColor = 6
$ws.Range("B4").Interior.Color = 65535

# Move the selection the way the saved workbook shows it.
$ws.Range("A7:XFD18").Select()
